$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 78

# Columns A-D hold text that looks numeric/date-like ("2025-02-19", "07"),
# so force them to Text before assigning, then clear the number format again
# so the new row matches the un-styled text cells used elsewhere in the sheet.
$textRange = $ws.Range("A${row}:D${row}")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-02-19"
$ws.Cells.Item($row, 2).Value = "22:41:42"
$ws.Cells.Item($row, 3).Value = "Wednesday"
$ws.Cells.Item($row, 4).Value = "07"

$textRange.ClearFormats()

$ws.Cells.Item($row, 5).Value = 129594
$ws.Cells.Item($row, 6).Value = 140480
$ws.Cells.Item($row, 7).Value = 171640
$ws.Cells.Item($row, 8).Value = 160054
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 145999
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192478
$ws.Cells.Item($row, 14).Value = 114990
$ws.Cells.Item($row, 15).Value = 45851
$ws.Cells.Item($row, 16).Value = 29099
$ws.Cells.Item($row, 17).Value = 67438
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 47614
$ws.Cells.Item($row, 20).Value = -1
